# Daily attendance processing - 2025-11-22 09:42:33
# Normalizes the "Recorded By" (column G) values: when the literal "System"
# appears as the first entry of a comma-separated recorder list, move it to
# the end of the list (e.g. "System, foo@bar.com" -> "foo@bar.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "System,*") { continue }

    $parts = $val -split ", "
    if ($parts[0] -ne "System") { continue }

    $rest = $parts[1..($parts.Length - 1)]
    if ($rest.Length -eq 1 -and $rest[0] -eq "admin@admin.com") { continue }

    $newVal = ($rest -join ", ") + ", System"
    $cell.Value = $newVal
}
